# plantilla_carga_masiva.xlsx - add a "Correo" column with mailto: hyperlinks
# for the two existing rows (Administrador / Ingenieria contacts).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personal")

# New header
$ws.Range("E1").Value = "Correo"

# New data rows - setting the value first (so the hyperlink doesn't overwrite
# it with a generic "display" string) then attaching the mailto: hyperlink.
$ws.Range("E2").Value = "admin@mach.com"
$ws.Range("E3").Value = "ingeniero@mach.com"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:admin@mach.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:ingeniero@mach.com") | Out-Null

# Match the column width used in the authored template as closely as the
# host's column-width rounding allows.
$ws.Columns.Item(5).ColumnWidth = 29.7

# Leave the selection where the author left it before saving.
$ws.Range("E12").Select() | Out-Null
